{"js": "// Apply the stat-table edits described by the diff.\n// The document contains a single 46-row, 1-column table where each row\n// holds one benchmark statistic (as plain text in the row's single cell).\n// Several rows get their value replaced outright, and the final three\n// \"raw dump\" rows (which held a tab-separated list of values packed into\n// one run) are collapsed down to a single short value each.\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// row index -> new text value for that row's single cell\nconst edits = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"808\",\n  5: \"0.00353\",\n  6: \"0.00017\",\n  7: \"0.00010\",\n  8: \"0.00024\",\n  9: \"0.00025\",\n  10: \"0.00035\",\n  11: \"0.14121\",\n  43: \"99.95\",\n  44: \"0.14\",\n  45: \"276\",\n};\n\nfor (const [rowIndex, newText] of Object.entries(edits)) {\n  const cell = table.getCell(Number(rowIndex), 0);\n  cell.value = newText;\n}\n\nawait context.sync();\n", "ps1": "# Apply the stat-table edits described by the diff.\n# The document contains a single 46-row, 1-column table where each row\n# holds one benchmark statistic (as plain text in the row's single cell).\n# Several rows get their value replaced outright, and the final three\n# \"raw dump\" rows (which held a tab-separated list of values packed into\n# one run) are collapsed down to a single short value each.\n#\n# Word COM table cells are 1-indexed, so table row N (0-indexed, as seen\n# in the XML / Office.js world) is Cell(N + 1, 1) here.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$edits = @{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"808\"\n    6  = \"0.00353\"\n    7  = \"0.00017\"\n    8  = \"0.00010\"\n    9  = \"0.00024\"\n    10 = \"0.00025\"\n    11 = \"0.00035\"\n    12 = \"0.14121\"\n    44 = \"99.95\"\n    45 = \"0.14\"\n    46 = \"276\"\n}\n\nforeach ($rowNumber in $edits.Keys) {\n    $t.Cell($rowNumber, 1).Range.Text = $edits[$rowNumber]\n}\n"}
